$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("B7").Value = -0.00789473684210527
$ws.Range("C7").Value = 0.7184210526315791
$ws.Range("D7").Value = 0.9828947368421049
$ws.Range("E7").Value = 0.9914104784810905
$ws.Range("F7").Value = 1.004686741267973
$ws.Range("G7").Value = 38

# Row 8
$ws.Range("B8").Value = 0.01891891891891893
$ws.Range("C8").Value = 0.7432432432432432
$ws.Range("D8").Value = 0.8883783783783785
$ws.Range("E8").Value = 0.9425382636149996
$ws.Range("F8").Value = 0.9553468939017323
$ws.Range("G8").Value = 37

# Row 9
$ws.Range("B9").Value = -0.08499999999999999
$ws.Range("C9").Value = 0.725
$ws.Range("D9").Value = 1.0145
$ws.Range("E9").Value = 1.007223907579641
$ws.Range("F9").Value = 1.02970358535076
$ws.Range("G9").Value = 20

# Row 10
$ws.Range("B10").Value = -0.08461538461538459
$ws.Range("C10").Value = 0.7769230769230769
$ws.Range("D10").Value = 1.114615384615385
$ws.Range("E10").Value = 1.055753467725958
$ws.Range("F10").Value = 1.095328074023299
$ws.Range("G10").Value = 13

# Row 11
$ws.Range("B11").Value = -0.26
$ws.Range("C11").Value = 0.42
$ws.Range("D11").Value = 0.3139999999999999
$ws.Range("E11").Value = 0.5603570290448759
$ws.Range("F11").Value = 0.5549774770204643
$ws.Range("G11").Value = 5
